# Fix a typo-era run-split in the last slide's "TextBox 5" paragraph.
# The paragraph's visible text is unchanged; three adjacent runs
# ("... In most cases ", "the system ", "is constructed ... makes ")
# are merged back into a single run so they share one consistent
# run-properties element (sz=2000, dirty=0).

$p  = $ppt.ActivePresentation
$s  = $p.Slides.Item(5)
$sh = $s.Shapes.Item(5)
$tr = $sh.TextFrame.TextRange

# "GenevaERS" occupies characters 1-9, so the run to normalize starts
# at character 10 and spans the combined length of the three runs
# being merged (95 + 11 + 79 = 185 characters).
$seg = $tr.Characters(10, 185)
$seg.Text = " is designed to resolve all processes in a single pass of the transaction file.  In most cases the system is constructed with pre-determined partitioning schemes.  Thus GenevaERS makes "
